$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data is price/volume text that looks numeric (e.g. "1.00",
# "0.0690", "43.023.58"). Excel auto-converts such strings when assigned
# via .Value, so force each touched cell/row range to Text format first.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "43.023.58"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "2.298.56"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "99.42"
$ws.Range("E6").Value = "  +1.53%  "

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "0.506"
$ws.Range("E7").Value = "  -0.66%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "36.28"
$ws.Range("E10").Value = "  +7.14%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.59%  "

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "17.98"
$ws.Range("E13").Value = "  +4.32%  "

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "6.92"
$ws.Range("E14").Value = "  +1.88%  "

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "2.655.67"
$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "2.295.86"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "0.798"
$ws.Range("E17").Value = "  -2.30%  "

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "42.947.48"
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "12.71"
$ws.Range("E19").Value = "  +8.71%  "

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "67.95"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "235.74"
$ws.Range("E23").Value = "  -0.60%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.83%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.46%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "24.95"
$ws.Range("E27").Value = "  +1.72%  "

$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = "168.73"
$ws.Range("E28").Value = "  +0.68%  "

$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D29").Value = "34.45"
$ws.Range("E29").Value = "  +0.43%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.87%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("D33").Value = "5.04"
$ws.Range("E33").Value = "  +1.17%  "

$ws.Range("B34:E34").NumberFormat = "@"
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "17.64"
$ws.Range("E34").Value = "  +4.45%  "

$ws.Range("B35:E35").NumberFormat = "@"
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "4.65"
$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.08%  "

$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0690"
$ws.Range("E37").Value = "  -1.63%  "

$ws.Range("B38:E38").NumberFormat = "@"
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.101"
$ws.Range("E38").Value = "  -0.20%  "

$ws.Range("B39:E39").NumberFormat = "@"
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "2.82"
$ws.Range("E39").Value = "  -0.83%  "

$ws.Range("B40:E40").NumberFormat = "@"
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.79"
$ws.Range("E40").Value = "  +0.66%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.00%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.85%  "

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "1.978.90"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "10.16"
$ws.Range("E45").Value = "  +2.54%  "

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "2.92"
$ws.Range("E46").Value = "  +1.63%  "

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "17.52"
$ws.Range("E47").Value = "  -1.27%  "

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "55.55"
$ws.Range("E48").Value = "  +3.72%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.62%  "

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "2.522.21"
$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "70.71"
$ws.Range("E51").Value = "  +0.47%  "

